# Implements AliExpressSite and first draft of Page Objects
#
# The only real content change in this revision is a copy/paste fix in the
# "PopUoPage" Page-Object notes (cell E7 on the "add song" sheet): the typo
# "PopUoPage" is corrected to "PopUpPage" and the stray blank line before the
# "SearchPage:" section is removed. Saving this change naturally reshuffles
# the shared-string table the same way Excel does (the old string is dropped
# and the corrected text is appended as a new shared string), which is why
# the sheet's cell/value indices shift even though every other cell's text is
# unchanged.
#
# The sheet view was also scrolled/zoomed in a bit further by the author
# (zoom 110% -> 140%, selection moved from E11 to E10) while reviewing the
# fix, so we reproduce that too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("add song")

$fixedText = "PopUpPage:`nclosePopUpIcon: xpath=""//div[contains(@class, 'dialog')]//a[.='x']""`npopUpHeader: xpath=""//div[contains(@class, 'dialog')]//h4""`nSearchPage:`nsearchTxt: id=""search-key""`nsearchBtn: css="".search-button"""

$ws.Range("E7").Value = $fixedText

# Match the author's end-state view: zoomed to 140%, scrolled so column D is
# the leftmost visible column, row 2 at top, and E10 selected/active.
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 140
$win.ScrollRow = 2
$win.ScrollColumn = 4
[void]$ws.Range("E10").Select()
